$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Add the new "route_protocol" (ospf) value and "process" column ---
$ws.Range("E2").Value = "ospf"
$ws.Range("F1").Value = "process"

# Fill in the remaining "ospf" / "1" values for the existing rows (2-7)
$ws.Range("E3").Value = "ospf"
$ws.Range("E4").Value = "ospf"
$ws.Range("E5").Value = "ospf"
$ws.Range("E6").Value = "ospf"
$ws.Range("E7").Value = "ospf"

$ws.Range("F2").Value = 1
$ws.Range("F3").Value = 1
$ws.Range("F4").Value = 1
$ws.Range("F5").Value = 1
$ws.Range("F6").Value = 1
$ws.Range("F7").Value = 1

# --- Add the two new devices (192.168.0.213 / 192.168.0.214), each with
# three interfaces / routes, mirroring the layout of the existing rows. ---

# Device names first (rows 8 and 11 start each new device block)
$ws.Range("A8").Value = "192.168.0.213"
$ws.Range("A11").Value = "192.168.0.214"

# Fill in the rest of column A for the new device blocks
$ws.Range("A9").Value = "192.168.0.213"
$ws.Range("A10").Value = "192.168.0.213"
$ws.Range("A12").Value = "192.168.0.214"
$ws.Range("A13").Value = "192.168.0.214"

# Column B (interfaces) for the new rows
$ws.Range("B8").Value = "GigabitEthernet2"
$ws.Range("B9").Value = "GigabitEthernet3"
$ws.Range("B10").Value = "GigabitEthernet4"
$ws.Range("B11").Value = "GigabitEthernet2"
$ws.Range("B12").Value = "GigabitEthernet3"
$ws.Range("B13").Value = "GigabitEthernet4"

# Column C (ip_address) for the new rows, in row order
$ws.Range("C8").Value = "10.100.34.3"
$ws.Range("C9").Value = "10.100.13.3"
$ws.Range("C10").Value = "10.100.23.3"
$ws.Range("C11").Value = "10.100.34.4"
$ws.Range("C12").Value = "10.100.24.4"
$ws.Range("C13").Value = "10.100.14.4"

# Column D (subnet_mask) for the new rows
$ws.Range("D8").Value = "255.255.255.0"
$ws.Range("D9").Value = "255.255.255.0"
$ws.Range("D10").Value = "255.255.255.0"
$ws.Range("D11").Value = "255.255.255.0"
$ws.Range("D12").Value = "255.255.255.0"
$ws.Range("D13").Value = "255.255.255.0"

# Column E (route_protocol) for the new rows
$ws.Range("E8").Value = "ospf"
$ws.Range("E9").Value = "ospf"
$ws.Range("E10").Value = "ospf"
$ws.Range("E11").Value = "ospf"
$ws.Range("E12").Value = "ospf"
$ws.Range("E13").Value = "ospf"

# Column F (process) for the new rows
$ws.Range("F8").Value = 1
$ws.Range("F9").Value = 1
$ws.Range("F10").Value = 1
$ws.Range("F11").Value = 1
$ws.Range("F12").Value = 1
$ws.Range("F13").Value = 1

# Move / update the current selection to match the saved view (C10)
$ws.Range("C10").Select()
